$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("setUp")
$ws2 = $wb.Worksheets.Item("Tests")

# --- setUp sheet: insert a new "pos" column (E) before the old F (test-results) block ---
$ws1.Columns("E:E").Insert()
$ws1.Columns("E:E").ColumnWidth = 5.3

# New "pos" field on the car datatype (row 8) and CarData index table (row 12 header + rows 14-23 values)
$ws1.Range("B8").Value = "int"
$ws1.Range("C8").Value = "pos"
$ws1.Range("E12").Value = "pos"
$ws1.Range("E14").Value = 1
$ws1.Range("E15").Value = 2
$ws1.Range("E16").Value = 3
$ws1.Range("E17").Value = 4
$ws1.Range("E18").Value = 5
$ws1.Range("E19").Value = 6
$ws1.Range("E20").Value = 7
$ws1.Range("E21").Value = 8
$ws1.Range("E22").Value = 9
$ws1.Range("E23").Value = 10

# Extend the "Data car CarData" header merge to cover the new column
$ws1.Range("B11:D11").UnMerge()
$ws1.Range("B11:E11").Merge()

# New test blocks: transformToWithPredicate2 / transformToUniquePredicate2
$ws1.Range("G9").Copy()
$ws1.Range("G13").PasteSpecial(-4122)
$ws1.Range("G13").Value = "Method Integer[] transformToWithPredicate2(car[] arr)"

$ws1.Range("G10").Copy()
$ws1.Range("G14").PasteSpecial(-4122)
$ws1.Range("G14").Value = "return arr[(x) transform to isEmpty(x) ? -1 : pos];"

$ws1.Range("G24").Copy()
$ws1.Range("G28").PasteSpecial(-4122)
$ws1.Range("G28").Value = "Method Integer[] transformToUniquePredicate2(car[] arr)"

$ws1.Range("G25").Copy()
$ws1.Range("G29").PasteSpecial(-4122)
$ws1.Range("G29").Value = "return arr[(x) transform unique to isEmpty(x) ? -1 : pos ];"

# --- Sheet view / selection housekeeping: setUp becomes the active tab ---
$ws2.Range("I14").Select()
$ws1.Activate()
$ws1.Range("I30").Select()
